# Daily attendance processing - reorder "Recorded By" names in column G.
# The exact text of several distinct "Recorded By" values is rewritten by
# rotating the list of names/emails (moving the last item to the front).
# This affects every cell in column G that holds one of these exact values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$col = $ws.Columns.Item(7)  # Column G - "Recorded By"

$col.Replace("backup@backdoor.com, System, system", "system, backup@backdoor.com, System", 1)
$col.Replace("dnasr281@gmail.com, System", "System, dnasr281@gmail.com", 1)
$col.Replace("System, admin@admin.com", "admin@admin.com, System", 1)
$col.Replace("dnasr281@gmail.com, admin@admin.com", "admin@admin.com, dnasr281@gmail.com", 1)
